$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of the data block (rows 2:3),
# pushing all existing data rows down by two (old row 2 -> new row 4, etc.)
$ws.Rows("2:3").Insert()

# The insert copies formatting from the row above (the bold/bordered header),
# which is not what the source rows had. Reset formatting on the new rows to
# the plain (unstyled) look used by the rest of the data rows.
$ws.Range("A2:T3").ClearFormats()

# Column D (Fecha) is a date column elsewhere formatted as
# "YYYY-MM-DD HH:MM:SS" - restore that number format for the new cells.
$ws.Range("D2:D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 2: Castle Brite / Primera
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = "Macroferia Regional de Talca"
$ws.Range("C2").Value = "Maule"
$ws.Range("D2").Value = "2021-12-15"
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100103
$ws.Range("H2").Value = "Frutos de hueso (carozo)"
$ws.Range("I2").Value = 100103003
$ws.Range("J2").Value = "Damasco"
$ws.Range("K2").Value = "Castle Brite"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 120
$ws.Range("N2").Value = 17000
$ws.Range("O2").Value = 17000
$ws.Range("P2").Value = 17000
$ws.Range("Q2").Value = "$/caja 15 kilos"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1133
$ws.Range("T2").Value = 15

# Row 3: Castle Brite / Segunda
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = "Macroferia Regional de Talca"
$ws.Range("C3").Value = "Maule"
$ws.Range("D3").Value = "2021-12-15"
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100103
$ws.Range("H3").Value = "Frutos de hueso (carozo)"
$ws.Range("I3").Value = 100103003
$ws.Range("J3").Value = "Damasco"
$ws.Range("K3").Value = "Castle Brite"
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = "$/caja 15 kilos"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1000
$ws.Range("T3").Value = 15
